# Set column O (CDS) values to "None" for all data rows (2 through 329),
# per fix: number of records/genes should be fetched from features and
# genes tables instead of the prior computed CDS summary text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 15).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 329 }

$ws.Range("O2:O$lastRow").Value = "None"
